$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp in the title cell (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 16:21"

# --- Update numeric stats for countries whose data changed ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2892096
$ws.Range("C4").Value = 1508
$ws.Range("D4").Value = 1235995
$ws.Range("E4").Value = 1523985
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 132116

# Row 7: India
$ws.Range("B7").Value = 651315
$ws.Range("C7").Value = 1426
$ws.Range("D7").Value = 395578
$ws.Range("E7").Value = 237042

# Row 18: Alemania
$ws.Range("B18").Value = 197210
$ws.Range("C18").Value = 210
$ws.Range("E18").Value = 6836
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 9074

# Row 26: Argentina
$ws.Range("D26").Value = 25930
$ws.Range("E26").Value = 45403
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 1453

# Row 29: Bielorrusia
$ws.Range("B29").Value = 63270
$ws.Range("C29").Value = 273
$ws.Range("D29").Value = 50669
$ws.Range("E29").Value = 12183
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 418

# Row 31: Belgica
$ws.Range("B31").Value = 61838
$ws.Range("C31").Value = 111
$ws.Range("D31").Value = 17091
$ws.Range("E31").Value = 34976
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = 9771

# Row 36: Kuwait
$ws.Range("B36").Value = 49303
$ws.Range("C36").Value = 631
$ws.Range("D36").Value = 39943
$ws.Range("E36").Value = 8995
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 365

# Row 41: Portugal
$ws.Range("B41").Value = 43569
$ws.Range("C41").Value = 413
$ws.Range("D41").Value = 28772
$ws.Range("E41").Value = 13192
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 1605

# Row 63: Nepal
$ws.Range("B63").Value = 15491
$ws.Range("C63").Value = 232
$ws.Range("D63").Value = 6415
$ws.Range("E63").Value = 9042
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 34

# Row 140: Uganda
$ws.Range("D140").Value = 868
$ws.Range("E140").Value = 59

# Rows 158-159: Vietnam / Namibia swap ranking (Namibia now ahead of Vietnam)
$ws.Range("A158").Value = "Namibia"
$ws.Range("B158").Value = 375
$ws.Range("C158").Value = 25
$ws.Range("D158").Value = 25
$ws.Range("E158").Value = 350

$ws.Range("A159").Value = "Vietnam"
$ws.Range("B159").Value = 355
$ws.Range("D159").Value = 340
$ws.Range("E159").Value = 15

# Rows 161-163: Isla de Man / Angola / Siria rotate ranking
$ws.Range("A161").Value = "Siria"
$ws.Range("B161").Value = 338
$ws.Range("C161").Value = 10
$ws.Range("D161").Value = 123
$ws.Range("E161").Value = 205
$ws.Range("H161").Value = 10

$ws.Range("A162").Value = "Isla de Man"
$ws.Range("B162").Value = 336
$ws.Range("D162").Value = 312
$ws.Range("E162").Value = 0
$ws.Range("H162").Value = 24

$ws.Range("A163").Value = "Angola"
$ws.Range("D163").Value = 107
$ws.Range("E163").Value = 203
$ws.Range("H163").Value = 18

# Rows 205-206: Dominica / Fiyi swap ranking (tied totals, only names swap)
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"
